# Fix the academic-pdf-report bug:
#  - D1 header should read "Personal Website" instead of "Google Personal Website"
#  - Restore the sheet's active selection to G3 (was erroneously left at G4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Personal Website"

$ws.Range("G3").Select()
